# "Add files via upload" -- append the newly played "pelada" entries to the
# Jogadores sheet (rows 443-463) and bring the filter / view state up to
# date with the larger data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jogadores")

# Each entry: Player, Vitorias, Empate, Derrotas, Gols, Partidas,
#             TardeDeVitoria, LaBarca, CraqueDoDia, GolsSofridos
$newRows = @(
    @("Coxinha", 4, 3, 2, 2, 1, 1, 0, 0, 0),
    @("Davi", 4, 3, 2, 1, 1, 1, 0, 0, 0),
    @("Leandrinho", 4, 3, 2, 3, 1, 1, 0, 0, 0),
    @("Jorge", 4, 3, 2, 3, 1, 1, 0, 1, 0),
    @("Leandrão", 4, 3, 2, 1, 1, 1, 0, 0, 0),
    @("Michel", 3, 1, 4, 0, 1, 0, 0, 0, 0),
    @("Corinthiano", 3, 1, 4, 2, 1, 0, 0, 0, 0),
    @("Juscielio", 3, 1, 4, 4, 1, 0, 0, 0, 0),
    @("Miqueias", 3, 1, 4, 1, 1, 0, 0, 0, 0),
    @("Ismael", 3, 1, 4, 1, 1, 0, 0, 0, 0),
    @("Digão", 2, 2, 3, 1, 1, 0, 1, 0, 0),
    @("Ranyeri", 2, 2, 3, 2, 1, 0, 1, 0, 0),
    @("Adriano", 2, 2, 3, 2, 1, 0, 1, 0, 0),
    @("Marcos", 2, 2, 3, 0, 1, 0, 1, 0, 0),
    @("Guinha", 2, 2, 3, 0, 1, 0, 1, 0, 0),
    @("Athos", 4, 0, 4, 0, 1, 0, 0, 0, 0),
    @("Du", 4, 0, 4, 0, 1, 0, 0, 0, 0),
    @("Marcelão", 4, 0, 4, 1, 1, 0, 0, 0, 0),
    @("Cabeleira", 4, 0, 4, 4, 1, 0, 0, 0, 0),
    @("Boneco", 4, 0, 4, 2, 1, 0, 0, 0, 0),
    @("Matheus", 9, 3, 4, 0, 1, 1, 0, 0, 12)
)

$startRow = 443
$filterLastRow = 446      # autofilter / _FilterDatabase only grew to here

function Write-PlayerRow($r, $entry) {
    $ws.Cells.Item($r, 1).Value = $entry[0]   # A - Jogadores
    $ws.Cells.Item($r, 3).Value = $entry[1]   # C - Vitorias
    $ws.Cells.Item($r, 4).Value = $entry[2]   # D - Empate
    $ws.Cells.Item($r, 5).Value = $entry[3]   # E - Derrotas
    $ws.Cells.Item($r, 6).Value = $entry[4]   # F - Gols
    $ws.Cells.Item($r, 7).Value = $entry[5]   # G - Partidas
    $ws.Cells.Item($r, 8).Value = $entry[6]   # H - Tarde de Vitoria
    $ws.Cells.Item($r, 9).Value = $entry[7]   # I - La barca
    $ws.Cells.Item($r, 10).Value = $entry[8]  # J - Craque do Dia
    $ws.Cells.Item($r, 11).Value = $entry[9]  # K - Gols Sofridos
}

# Write the rows up through the old filter's new edge first ...
$row = $startRow
for ($i = 0; $i -lt $newRows.Count -and $row -le $filterLastRow; $i++) {
    Write-PlayerRow $row $newRows[$i]
    $row = $row + 1
}

# ... then grow the autofilter to match (toggle off/on since it only
# resizes when re-applied), and keep the _FilterDatabase name matching it.
$ws.Range("A1:K421").AutoFilter() | Out-Null
$ws.Range("A1:K$filterLastRow").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Jogadores!`$A`$1:`$K`$$filterLastRow"
    }
}

# ... and only now add the remaining rows, which fall outside the filter
# range (matches the source workbook, where the filter wasn't re-applied
# after the last few rows were pasted in).
for ($i; $i -lt $newRows.Count; $i++) {
    Write-PlayerRow $row $newRows[$i]
    $row = $row + 1
}

$lastRow = $row - 1   # 463 - last row actually written

# Leave the selection where Excel would after typing the last new row in
# by hand -- one column to the right of the data, on the last data row.
$ws.Range("L$lastRow").Select()
